# MODFLOW 6 GWT planning.xlsx
# Worked on times series of MAW, CSUB, SFR, LAK, and UZF packages.
# Worked on supporting times series in the formula editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 "Time Lists": TIME series used count bumped from 50 to 51
$ws.Range("I3").Value = 51

# "TIME series used" column (I) gained entries for packages worked on
# Row 24: CSUB
$ws.Range("I24").Value = 0.5
# Row 27: LAK
$ws.Range("I27").Value = 1
# Row 28: SFR
$ws.Range("I28").Value = 1
# Row 29: MAW
$ws.Range("I29").Value = 0.75
# Row 35: UZF
$ws.Range("I35").Value = 2

# Totals block (rows 39-42) extended into column I to mirror column H
$ws.Range("I39").Formula = "=SUM(I2:I37)"
$ws.Range("I40").Formula = "=I39/5"
$ws.Range("I41").Formula = "=I40*7/5"
$ws.Range("I42").Formula = "=I41/30"

# Leave the selection on the cell that was last edited
$ws.Range("I3").Select() | Out-Null
